# Auto-applies the crypto price/volume update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.995.65"
$ws.Range("E2").Value = "  +5.89%  "
$ws.Range("D3").Value = "2.331.80"
$ws.Range("E3").Value = "  +5.00%  "
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.59"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.11"
$ws.Range("E6").Value = "  +9.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  +4.39%  "
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  +9.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.10"
$ws.Range("E10").Value = "  +7.39%  "
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  +8.09%  "
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "2.690.63"
$ws.Range("E14").Value = "  +4.98%  "
$ws.Range("D15").Value = "2.335.44"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.841"
$ws.Range("E16").Value = "  +4.83%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.19"
$ws.Range("E17").Value = "  +7.87%  "
$ws.Range("D18").Value = "46.879.66"
$ws.Range("E18").Value = "  +6.16%  "
$ws.Range("E19").Value = "  +20.96%  "
$ws.Range("D20").Value = "0.0₃0957"
$ws.Range("E20").Value = "  +5.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.22"
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.29"
$ws.Range("E22").Value = "  +5.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.51"
$ws.Range("E23").Value = "  +9.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  +3.99%  "
$ws.Range("E25").Value = "  +5.30%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.27"
$ws.Range("E27").Value = "  +17.00%  "
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.97"
$ws.Range("E29").Value = "  +5.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.41"
$ws.Range("E30").Value = "  +4.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.86"
$ws.Range("E31").Value = "  +4.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0823"
$ws.Range("E32").Value = "  +9.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.07"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.65"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.115"
$ws.Range("E35").Value = "  +9.05%  "
$ws.Range("E36").Value = "  +6.16%  "
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("E38").Value = "  +3.87%  "
$ws.Range("E39").Value = "  +10.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0312"
$ws.Range("E40").Value = "  +9.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  +6.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.08"
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  +16.58%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.63"
$ws.Range("E45").Value = "  +17.46%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.802.98"
$ws.Range("E46").Value = "  +3.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.98"
$ws.Range("E47").Value = "  +12.17%  "
$ws.Range("E48").Value = "  +8.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.17"
$ws.Range("E49").Value = "  +4.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.42"
$ws.Range("E50").Value = "  +6.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.84"
$ws.Range("E51").Value = "  +3.03%  "
